# Updated symbol list / crypto price refresh
# ------------------------------------------------------------------
# NOTE on this engine's quirks (discovered empirically):
#  * A worksheet variable captured via $wb.ActiveSheet / Worksheets.Item()
#    can silently start pointing at a DIFFERENT sheet once another sheet
#    is added/activated/deleted. To stay safe we always re-resolve the
#    worksheet via $wb.Worksheets.Item($wsName) immediately before each
#    use, rather than trusting a long-lived $ws variable.
#  * Range.Value, when assigned a numeric-looking string (e.g. "268.85"),
#    gets silently coerced into a real number - which would not match the
#    original text/string cells in this workbook. To keep such values as
#    TEXT we stage them in a helper cell formatted as Text ("@") on a
#    throwaway worksheet, then use Copy + PasteSpecial(xlPasteValues) to
#    move the text value onto the destination cell without altering the
#    destination's own style/formatting.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsName = $wb.ActiveSheet.Name

# --- 1) Plain text cell updates (coin names / links / rank labels) ---
# These are not numeric-looking, so direct assignment keeps them as text.
$wb.Worksheets.Item($wsName).Range("B15").Value = "ProBitToken"
$wb.Worksheets.Item($wsName).Range("C15").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$wb.Worksheets.Item($wsName).Range("E15").Value = "14ProBitTokenPROB"
$wb.Worksheets.Item($wsName).Range("B16").Value = "BitMartToken"
$wb.Worksheets.Item($wsName).Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$wb.Worksheets.Item($wsName).Range("E16").Value = "15BitMartTokenBMX"
$wb.Worksheets.Item($wsName).Range("B17").Value = "MCDex"
$wb.Worksheets.Item($wsName).Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$wb.Worksheets.Item($wsName).Range("E17").Value = "16MCDexMCB"
$wb.Worksheets.Item($wsName).Range("B18").Value = "BitForexToken"
$wb.Worksheets.Item($wsName).Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$wb.Worksheets.Item($wsName).Range("E18").Value = "17BitForexTokenBF"
$wb.Worksheets.Item($wsName).Range("B19").Value = "CoinExToken"
$wb.Worksheets.Item($wsName).Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$wb.Worksheets.Item($wsName).Range("E19").Value = "18CoinExTokenCET"
$wb.Worksheets.Item($wsName).Range("B20").Value = "TigerCash"
$wb.Worksheets.Item($wsName).Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$wb.Worksheets.Item($wsName).Range("E20").Value = "19TigerCashTCH"
$wb.Worksheets.Item($wsName).Range("B21").Value = "HotbitToken"
$wb.Worksheets.Item($wsName).Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$wb.Worksheets.Item($wsName).Range("E21").Value = "20HotbitTokenHTB"
$wb.Worksheets.Item($wsName).Range("B22").Value = "BitKan"
$wb.Worksheets.Item($wsName).Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$wb.Worksheets.Item($wsName).Range("E22").Value = "21BitKanKAN"
$wb.Worksheets.Item($wsName).Range("B23").Value = "NitroEx"
$wb.Worksheets.Item($wsName).Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$wb.Worksheets.Item($wsName).Range("E23").Value = "22NitroExNTX"
$wb.Worksheets.Item($wsName).Range("B24").Value = "LEO"
$wb.Worksheets.Item($wsName).Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$wb.Worksheets.Item($wsName).Range("E24").Value = "23LEOLEO"
$wb.Worksheets.Item($wsName).Range("B25").Value = "BTSEToken"
$wb.Worksheets.Item($wsName).Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$wb.Worksheets.Item($wsName).Range("E25").Value = "24BTSETokenBTSE"
$wb.Worksheets.Item($wsName).Range("B26").Value = "BitpandaEcosystemToken"
$wb.Worksheets.Item($wsName).Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$wb.Worksheets.Item($wsName).Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$wb.Worksheets.Item($wsName).Range("B42").Value = "CEJI"
$wb.Worksheets.Item($wsName).Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$wb.Worksheets.Item($wsName).Range("E42").Value = "41CEJICEJIWorstin24h"
$wb.Worksheets.Item($wsName).Range("B43").Value = "BKEXToken"
$wb.Worksheets.Item($wsName).Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$wb.Worksheets.Item($wsName).Range("E43").Value = "42BKEXTokenBKK"
$wb.Worksheets.Item($wsName).Range("B49").Value = "BOLO"
$wb.Worksheets.Item($wsName).Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$wb.Worksheets.Item($wsName).Range("E49").Value = "48BOLOBOLO"
$wb.Worksheets.Item($wsName).Range("B50").Value = "CryptobidCoin"
$wb.Worksheets.Item($wsName).Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$wb.Worksheets.Item($wsName).Range("E50").Value = "49CryptobidCoinCBC"

# --- 2) Numeric-looking price cells, kept as TEXT via helper sheet ---
$tempSheet = $wb.Worksheets.Add()
$tempName = $tempSheet.Name
$wb.Worksheets.Item($tempName).Range("A1").NumberFormat = "@"

$wb.Worksheets.Item($tempName).Range("A1").Value = "268.85"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D2").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "22.88"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D3").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "6.327"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D4").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.06177"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D5").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "3.642"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D6").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "6.674"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D7").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "1.383"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D8").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.8293"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D9").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.01375"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D10").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.1603"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D11").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.08268"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D12").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.03485"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D13").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.03237"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D14").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.1239"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D15").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.09314"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D16").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "3.839"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D17").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.001647"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D18").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.04731"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D19").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.006339"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D20").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.005665"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D21").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.001076"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D22").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.0001499"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D23").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "3.721"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D24").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "2.413"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D25").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.3340"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D26").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.0002702"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D27").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.04688"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D40").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.006924"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D41").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.003797"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D42").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.1158"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D43").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.01153"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D44").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.00006246"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D45").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.0009898"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D46").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.00000000749"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D47").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.9193"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D48").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.002229"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D49").PasteSpecial(-4163)

$wb.Worksheets.Item($tempName).Range("A1").Value = "0.00001399"
$wb.Worksheets.Item($tempName).Range("A1").Copy()
$wb.Worksheets.Item($wsName).Range("D50").PasteSpecial(-4163)

# --- 3) Clean up the helper sheet and clipboard state ---
$excel.CutCopyMode = 0
$excel.DisplayAlerts = $false
$wb.Worksheets.Item($tempName).Delete()
$excel.DisplayAlerts = $true

$wb.Worksheets.Item($wsName).Activate()
